$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target confusion-matrix / metric values for rows 2-11 (columns C:I)
# C = tp, D = fp, E = tn, F = fn, G = precision, H = recall, I = fscore
$data = @{
    2  = @(135, 16, 0,  0,  0.8940397350993378, 1,                  0.9440559440559442)
    3  = @(135, 16, 0,  0,  0.8940397350993378, 1,                  0.9440559440559442)
    4  = @(135, 16, 0,  0,  0.8940397350993378, 1,                  0.9440559440559442)
    5  = @(135, 16, 0,  0,  0.8940397350993378, 1,                  0.9440559440559442)
    6  = @(135, 16, 0,  0,  0.8940397350993378, 1,                  0.9440559440559442)
    7  = @(135, 16, 0,  0,  0.8940397350993378, 1,                  0.9440559440559442)
    8  = @(132, 16, 0,  3,  0.8918918918918919, 0.9777777777777777, 0.9328621908127209)
    9  = @(43,  6,  10, 92, 0.8775510204081632, 0.3185185185185185, 0.4673913043478261)
    10 = @(43,  6,  10, 89, 0.8775510204081632, 0.3257575757575757, 0.4751381215469612)
    11 = @(41,  5,  11, 89, 0.8913043478260869, 0.3153846153846154, 0.4659090909090909)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 3).Value = $values[0]
    $ws.Cells.Item($row, 4).Value = $values[1]
    $ws.Cells.Item($row, 5).Value = $values[2]
    $ws.Cells.Item($row, 6).Value = $values[3]
    $ws.Cells.Item($row, 7).Value = $values[4]
    $ws.Cells.Item($row, 8).Value = $values[5]
    $ws.Cells.Item($row, 9).Value = $values[6]
}
